$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. "1.000", "0.9999") must be
# forced to Text format first, otherwise Excel will silently convert them to numeric
# values and strip formatting such as trailing zeros (e.g. "1.000" -> 1).
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '22.208.28'
$ws.Range('E2').Value = '  -1.25%  '
$ws.Range('D3').Value = '1.565.62'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '1.000'
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('D6').Value = '289.63'
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').Value = '0.3770'
$ws.Range('E7').Value = '  +2.33%  '
$ws.Range('D8').Value = '0.3297'
$ws.Range('E8').Value = '  -1.21%  '
$ws.Range('D9').Value = '44.43'
$ws.Range('E9').Value = '  -7.52%  '
$ws.Range('E10').Value = '  -0.06%  '
$ws.Range('D11').Value = '0.07391'
$ws.Range('E11').Value = '  -2.44%  '
$ws.Range('D12').Value = '1.000'
$ws.Range('E12').Value = '  -0.16%  '
$ws.Range('D13').Value = '20.17'
$ws.Range('E13').Value = '  -2.97%  '
$ws.Range('D14').Value = '5.868'
$ws.Range('E14').Value = '  -1.81%  '
$ws.Range('D15').Value = '6.919'
$ws.Range('E15').Value = '  -0.45%  '
$ws.Range('D16').Value = '1.572.15'
$ws.Range('E16').Value = '  -0.25%  '
$ws.Range('E17').Value = '  -2.70%  '
$ws.Range('D18').Value = '0.06642'
$ws.Range('E18').Value = '  -1.51%  '
$ws.Range('D19').Value = '85.84'
$ws.Range('E19').Value = '  -2.91%  '
$ws.Range('D20').Value = '6.442'
$ws.Range('E20').Value = '  +0.60%  '
$ws.Range('D21').Value = '0.9998'
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').Value = '16.19'
$ws.Range('E22').Value = '  -2.73%  '
$ws.Range('D23').Value = '11.78'
$ws.Range('E23').Value = '  -2.42%  '
$ws.Range('D24').Value = '22.205.67'
$ws.Range('E24').Value = '  -1.20%  '
$ws.Range('D25').Value = '2.272'
$ws.Range('E25').Value = '  -5.37%  '
$ws.Range('D26').Value = '2.550'
$ws.Range('E26').Value = '  -4.07%  '
$ws.Range('D27').Value = '150.45'
$ws.Range('E27').Value = '  -1.49%  '
$ws.Range('D28').Value = '19.18'
$ws.Range('E28').Value = '  -3.04%  '
$ws.Range('D29').Value = '4.920'
$ws.Range('E29').Value = '  -1.56%  '
$ws.Range('D30').Value = '1.744.37'
$ws.Range('E30').Value = '  -0.47%  '
$ws.Range('D31').Value = '122.06'
$ws.Range('E31').Value = '  -2.75%  '
$ws.Range('D32').Value = '1.120'
$ws.Range('E32').Value = '  +1.70%  '
$ws.Range('D33').Value = '6.086'
$ws.Range('E33').Value = '  -1.23%  '
$ws.Range('D34').Value = '1.881'
$ws.Range('E34').Value = '  -5.98%  '
$ws.Range('B35').Value = 'FraxShare'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D35').Value = '9.409'
$ws.Range('E35').Value = '  -4.65%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').Value = '0.08235'
$ws.Range('E36').Value = '  -1.98%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').Value = '5.328'
$ws.Range('E37').Value = '  -1.43%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.02330'
$ws.Range('E38').Value = '  -6.21%  '
$ws.Range('D39').Value = '0.06234'
$ws.Range('E39').Value = '  -3.55%  '
$ws.Range('D40').Value = '0.2152'
$ws.Range('E40').Value = '  -5.15%  '
$ws.Range('D41').Value = '1.244'
$ws.Range('E41').Value = '  -4.54%  '
$ws.Range('D42').Value = '11.14'
$ws.Range('E42').Value = '  -3.52%  '
$ws.Range('D43').Value = '0.6036'
$ws.Range('E43').Value = '  -4.77%  '
$ws.Range('D44').Value = '1.0000'
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('D45').Value = '13.80'
$ws.Range('E45').Value = '  -3.06%  '
$ws.Range('D46').Value = '3.759'
$ws.Range('E46').Value = '  -0.73%  '
$ws.Range('D47').Value = '0.5848'
$ws.Range('E47').Value = '  -5.41%  '
$ws.Range('D48').Value = '2.002'
$ws.Range('E48').Value = '  -3.85%  '
$ws.Range('D49').Value = '121.88'
$ws.Range('E49').Value = '  -3.51%  '
$ws.Range('D50').Value = '1.179'
$ws.Range('E50').Value = '  -2.93%  '
$ws.Range('D51').Value = '0.07019'
$ws.Range('E51').Value = '  -2.99%  '

# Restore the default cell style on the cells we forced to Text format, so the
# workbook does not end up with a stray/unused number format left on the cell.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}